$wb = $excel.ActiveWorkbook

# Helper: write a literal text label into a cell without Excel's automatic
# "looks like a number" coercion (a plain .Value = "2050" assignment would
# silently store 2050 as a number). Routing the text through a formula
# result and then freezing it with Copy/PasteSpecial(values) preserves the
# original cell style (no NumberFormat/quote-prefix side effects) while
# still landing as a genuine text cell.
function Set-TextLabel($range, [string]$text) {
    $escaped = $text -replace '"', '""'
    $range.Formula = "=""$escaped"""
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# Fix the mis-rendered "2050" column header (it was a leftover numeric
# value, 673.0764197344222, instead of the year/period label used by the
# other header cells in each table).
Set-TextLabel $wb.Worksheets.Item(1).Range("E1") "2050"
Set-TextLabel $wb.Worksheets.Item(2).Range("E1") "2050"
Set-TextLabel $wb.Worksheets.Item(3).Range("E1") "2050"
Set-TextLabel $wb.Worksheets.Item(4).Range("E1") "2041-2050"
Set-TextLabel $wb.Worksheets.Item(5).Range("E1") "2050"

# Remove the "Total" row from each table (row 13 on sheets 1-4, which
# spanned A1:E13; row 4 on sheet 6, which spanned A1:B4). Sheet 5 had no
# Total row to begin with.
foreach ($idx in 1..4) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Rows.Item(13).Delete()
}

$ws6 = $wb.Worksheets.Item(6)
$ws6.Rows.Item(4).Delete()
